# Edit workbook per target diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Update D33 text (remove 2 <img> lines, close <h4> tag) ----
$text109 = @'
Sensus Navigation(2016)이 구동 중인 경우 업데이트 방법에 대해 고객이 보고할 수 있습니다.
<br>
<br>차량이 네트워크가 동작하는 경우, 무선(OTA)로 업데이트할 수 있습니다.
<br><h3>인터넷을 통해 업데이트 방법</h3>
<br>1. 소프트웨어 업데이트는 운전 중이거나 차량의 시동이 켜져 있는 경우 차량에 직접 다운로드됩니다.
<br>2. 업데이트가 다운로드되면 업데이트를 설치할 준비가 되었다는 알림을 받게 됩니다. '지금 설치'를 클릭하여 설치를 시작합니다.
<br>즉시 설치하거나 원하는 시간을 선택하여 설치를 시작할 수 있습니다. 설치를 완료하는 데 약 90분이 소요됩니다.
<br>3.소프트웨가 최신 상태입니다. 메시지가 중앙 화면에 나타나면 소프트웨어 업데이트가 완료된 것입니다.
<br><h4>차량에서 직접 다운로드하기</h4>
<br>1. <img src="https://www.volvocars.com/images/support/imgc0eb4a9734129c4bc0a801522aaf58ce_2_--_--_VOICEpnghigh.png" width="64" height="64">
<br>앱 보기에서 다운로드 센터 버튼을 누릅니다.
<br>2. 지도 버튼을 누릅니다.
<br>3.설치 버튼을 누른 후 확인을 선택합니다.
<br>4. 선택한 지도 업데이트의 설치가 시작됩니다.
<br>
<br><h4>진행 중인 지도 다운로드가 취소될 경우</h4>
<br>지도 다운로드 중에 차량 시동을 끄면, 차량을 다시 시동하여 인터넷에 다시 연결할 때 지도 다운로드가 다시 시작됩니다.
<br><h3>지도를 USB로 설치하는 방법</h3>
<br>새 지도는 인터넷에 연결된 컴퓨터에서 USB 메모리로 다운로드한 후 USB 메모리에서 차량의 내비게이션 시스템으로 전송할 수 있습니다.
<br>
<br>USB 메모리
<br>USB 메모리로 업데이트를 관리해야 하는 경우에 다음 요구 사항이 적용됩니다.
<br>USB 표준:	최소 2.0
<br>파일 시스템: FAT32, exFAT 또는 NFTS
<br>용량: 최대 128GB
<br>
<br>지도 파일은 아래 링크에서 다운로드 받을 수 있습니다.
<br><img src="https://az685612.vo.msecnd.net/swdlimages/SPA_KR_v01_644x386.png">
<br><a href="https://www.volvocars.com/kr/support/downloads/maps/spa/daehanmingug" target="_blank">Sensus Navigation 2016 지도 다운로드</a>
<br>
<br><h3>현재 지도 버전 확인</h3>
<br>1. <img src="https://www.volvocars.com/images/support/imgc0eb4a9734129c4bc0a801522aaf58ce_2_--_--_VOICEpnghigh.png" width="64" height="64">
<br>앱 보기에서 다운로드 센터 버튼을 누릅니다.
<br>2. 지도 버튼을 누릅니다.
<br>> 대한민국을 누릅니다.
<br>3.<img="https://www.volvocars.com/images/support/imgf13a34931b25bc6ac0a8015256d334bc_1_--_--_VOICEpnghigh.png" width="64" height="64">
<br>화살표를 터치하여 선택한 지역을 확대합니다.
<br>> 상세 지도 정보가 표시됩니다.
<br>4. 지도 이미지 아래에서 지도 정보를 확인합니다.
'@
$ws.Cells.Item(33, 4).Value = $text109

# ---- Update D34 text (shrink img width/height from 64 to 20) ----
$text112 = @'
<br>지원 모델: iCup 지원 차량
<br>
<br>고객은 차량에 블루투스로 연결하면 음악이나 기타 미디어 사운드가 재생되지 않는다고 보고할 수 있습니다.
<br>홈 보기 또는 앱 보기<img src="https://www.volvocars.com/images/support/img7a8807df5cacedebc0a801525a46ef6f_2_--_--_VOICEpnghigh.png" width="20" height="20">에서 블루투스 미디어 플레이어 앱<img src="https://www.volvocars.com/images/support/img963609385e11fb97c0a801525c7b6fab_3_--_--_VOICEpnghigh.png" width="20" height="20">을 시작합니다.
<br> 앱 사용 시 홈 보기에서 단축 명령을 사용하여 앱을 제어할 수도 있습니다.
<br>
<br><em>참고
<br>Bluetooth를 통해 전화에서 미디어를 스트리밍하려면 먼저 블루투스 미디어 플레이어를 시작해야 합니다.</em>
<br>
<br>음성 컨트롤을 사용하여 미디어 플레이어를 제어할 수도 있습니다.
'@
$ws.Cells.Item(34, 4).Value = $text112

# ---- Add new row 35: "경보 레벨 낮추기" ----
$text115 = @'
고객이 차박을 하거나 캠핑을 할 때에 차량의 경보 장치로 인한 불편함을 보고할 수 있습니다.
<br>이러한 경우 아래 단계를 이용해 문제를 해결할 수 있습니다.
<br>
<br>경보 레벨 낮춤이 켜지면 경보장치의 동작 및 기울기 감지기가 꺼집니다. 그러면 경보장치가 차량 내의 움직임에 반응하지 않습니다.
<br>이 기능은 차량의 중앙 화면에서 또는 차량이 잠겨 있을 때 Volvo Cars 앱을 사용하여 켤 수 있습니다.
<br>
<br>경보 레벨 낮춤 기능은 각 사용 후 꺼지며 그러면 다시 켜야 합니다.
<br>
<br><h3>차량 중앙 화면에서 경보 레벨 감소 켜기</h3>
<br>1. 차량 중앙 화면에서 우측 하단의 설정을 누릅니다.
<br>2. 이후 컨트롤 메뉴를 선택한 후 알람 최소 모드를 켭니다.
<br>이 기능은 차량이 구동될 때까지 활성 상태이며, 그런 다음에 다시 켜야 합니다.
<br>경보 레벨 낮춤 기능은 설정 메뉴에서도 끌 수 있습니다.
'@
$ws.Cells.Item(35, 1).Value = "경보 레벨 낮추기"
$ws.Cells.Item(35, 2).Value = [DateTime]"2024-06-11"
$ws.Cells.Item(35, 3).Value = "HT203028"
$ws.Cells.Item(35, 4).Value = $text115
$ws.Cells.Item(35, 4).WrapText = $true
$ws.Rows.Item(35).RowHeight = 227.25

# ---- Add new row 36: "스티어링휠 잠금장치" ----
$ws.Cells.Item(36, 1).Value = "스티어링휠 잠금장치"
$ws.Cells.Item(36, 2).Value = [DateTime]"2024-06-11"
$ws.Cells.Item(36, 3).Value = "HT203029"

# ---- Update view: scroll to bottom, select D36 ----
$win = $excel.ActiveWindow
$win.ScrollRow = 34
$win.ScrollColumn = 1
$ws.Range("D36").Select() | Out-Null
